$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 email address (sk@g.com -> sj@g.com) - inserted before "remarks"
# to match shared-string table ordering in the target workbook.
$ws.Range("B5").Value = "sj@g.com"

# Add a new "remarks" column header in C1
$ws.Range("C1").Value = "remarks"

# Fill in remarks column based on email validity
$ws.Range("C2").Value = "Invalid"
$ws.Range("C3").Value = "Invalid"
$ws.Range("C4").Value = "Invalid"
$ws.Range("C5").Value = "Valid"

# Update selection to match target state
$ws.Range("D11").Select()
